$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Progression")

# Row 10 (45939) gets filled in with new entries. First mirror the
# cell formatting already used on row 9 (B=Neutral, C/E/G/H/I=Good,
# F/J=Neutral, D/K/L/M=Bad) so the same style entries get reused
# instead of creating new ones.
$ws.Range("B9:M9").Copy()
$ws.Range("B10:M10").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now set the values, in the order that matches how the new shared
# strings were appended to the workbook.
$ws.Range("B10").Value = "Outlander"
$ws.Range("C10").Value = "Puissance"
$ws.Range("G10").Value = "Coup critique"
$ws.Range("J10").Value = "Brigand"
$ws.Range("H10").Value = "Entraînement II"
$ws.Range("I10").Value = "Infusions magiques"
$ws.Range("E10").Value = "Cleptomane"
$ws.Range("F10").Value = "Mentor vénérable"

# Row 11 (45940) picks up an additional entry in column G, using the
# same "Good" formatting already present on the analogous cell G10/G9.
$ws.Range("G9").Copy()
$ws.Range("G11").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("G11").Value = "Entraînement II"
